# Update countries & provincias Spain
# Applies the daily COVID data refresh: a handful of per-country counters
# were updated, and (since the sheet is kept sorted by "Casos totales"
# descending) a few rows swapped places as a result. The timestamp caption
# in A1 is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp caption (row 1) -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 13:03"

# --- Plain per-country counter updates (no reordering) --------------------
# Sri Lanka (row 103): Casos activos, Recuperados
$ws.Cells.Item(103, 4).Value = 187
$ws.Cells.Item(103, 5).Value = 524

# Malta (row 117): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Cells.Item(117, 2).Value = 480
$ws.Cells.Item(117, 3).Value = 3
$ws.Cells.Item(117, 4).Value = 399
$ws.Cells.Item(117, 5).Value = 77

# Taiwan (row 121): Casos totales, Nuevos casos, Recuperados
$ws.Cells.Item(121, 2).Value = 438
$ws.Cells.Item(121, 3).Value = 6
$ws.Cells.Item(121, 5).Value = 98

# Madagascar (row 142): Casos activos, Recuperados
$ws.Cells.Item(142, 4).Value = 99
$ws.Cells.Item(142, 5).Value = 50

# --- Malaui gets fresh numbers and jumps above Mongolia / Puerto Rico -----
# Row 174 now holds Malaui's (updated) data, row 175 holds what used to be
# Mongolia's row, and row 176 holds what used to be Puerto Rico's row.
$ws.Cells.Item(174, 1).Value = "Malaui"
$ws.Cells.Item(174, 2).Value = 41
$ws.Cells.Item(174, 3).Value = 2
$ws.Cells.Item(174, 4).Value = 9
$ws.Cells.Item(174, 5).Value = 29
$ws.Cells.Item(174, 6).Value = 1
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 3

$ws.Cells.Item(175, 1).Value = "Mongolia"
$ws.Cells.Item(175, 2).Value = 40
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(175, 4).Value = 12
$ws.Cells.Item(175, 5).Value = 28
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

$ws.Cells.Item(176, 1).Value = "Puerto Rico"
$ws.Cells.Item(176, 2).Value = 39
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 1
$ws.Cells.Item(176, 5).Value = 36
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 2

# --- Belice swaps ahead of Santa Lucia (tied totals) -----------------------
$ws.Cells.Item(188, 1).Value = "Belice"
$ws.Cells.Item(188, 2).Value = 18
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 13
$ws.Cells.Item(188, 5).Value = 3
$ws.Cells.Item(188, 6).Value = 1
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 2

$ws.Cells.Item(189, 1).Value = "Santa Lucia"
$ws.Cells.Item(189, 2).Value = 18
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 15
$ws.Cells.Item(189, 5).Value = 3
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 0

# --- San Vicente y las Granadinas swaps ahead of Namibia (tied totals) ----
$ws.Cells.Item(194, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(194, 2).Value = 16
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 8
$ws.Cells.Item(194, 5).Value = 8
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

$ws.Cells.Item(195, 1).Value = "Namibia"
$ws.Cells.Item(195, 2).Value = 16
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 8
$ws.Cells.Item(195, 5).Value = 8
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

# --- San Cristobal y Nieves swaps ahead of Burundi (tied totals) ----------
$ws.Cells.Item(198, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(198, 2).Value = 15
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 8
$ws.Cells.Item(198, 5).Value = 7
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0

$ws.Cells.Item(199, 1).Value = "Burundi"
$ws.Cells.Item(199, 2).Value = 15
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 7
$ws.Cells.Item(199, 5).Value = 7
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 1
